$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29
$ws.Range("A29").Value = 130671357
$ws.Range("B29").Value = 79243
$ws.Range("E29").Value = 6425
$ws.Range("F29").Value = "Garnlav"
$ws.Range("G29").Value = "Alectoria sarmentosa"
$ws.Range("H29").Value = "(Ach.) Ach."
$ws.Range("Q29").Value = 557304
$ws.Range("R29").Value = 6710306
$ws.Range("Z29").Value = "09:53"
$ws.Range("AB29").Value = "09:53"

# Row 30
$ws.Range("A30").Value = 130671332
$ws.Range("B30").Value = 91808
$ws.Range("E30").Value = 1202
$ws.Range("F30").Value = "Ullticka"
$ws.Range("G30").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H30").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q30").Value = 556970
$ws.Range("R30").Value = 6710400
$ws.Range("Z30").Value = "11:06"
$ws.Range("AB30").Value = "11:06"

# Row 52
$ws.Range("A52").Value = 130671372
$ws.Range("B52").Value = 5177
$ws.Range("E52").Value = 100526
$ws.Range("F52").Value = "Bronshjon"
$ws.Range("G52").Value = "Callidium coriaceum"
$ws.Range("H52").Value = "Paykull, 1800"
$ws.Range("Q52").Value = 557068
$ws.Range("R52").Value = 6710403
$ws.Range("Z52").Value = "10:48"
$ws.Range("AB52").Value = "10:48"

# Row 53
$ws.Range("A53").Value = 130671325
$ws.Range("B53").Value = 5197
$ws.Range("D53").Value = "LC"
$ws.Range("E53").Value = 105930
$ws.Range("F53").Value = "Vågbandad barkbock"
$ws.Range("G53").Value = "Semanotus undatus"
$ws.Range("H53").Value = "(Linnaeus, 1758)"
$ws.Range("M53").Value = "färska gnagspår"
$ws.Range("Q53").Value = 557029
$ws.Range("R53").Value = 6710414
$ws.Range("Z53").Value = "10:51"
$ws.Range("AB53").Value = "10:51"

# Row 54
$ws.Range("A54").Value = 130671330
$ws.Range("B54").Value = 91808
$ws.Range("D54").Value = "NT"
$ws.Range("E54").Value = 1202
$ws.Range("F54").Value = "Ullticka"
$ws.Range("G54").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H54").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("M54").ClearContents()
$ws.Range("Q54").Value = 557102
$ws.Range("R54").Value = 6710446
$ws.Range("Z54").Value = "10:34"
$ws.Range("AB54").Value = "10:34"

# Row 60
$ws.Range("A60").Value = 130671326
$ws.Range("B60").Value = 5197
$ws.Range("D60").Value = "LC"
$ws.Range("E60").Value = 105930
$ws.Range("F60").Value = "Vågbandad barkbock"
$ws.Range("G60").Value = "Semanotus undatus"
$ws.Range("H60").Value = "(Linnaeus, 1758)"
$ws.Range("M60").Value = "färska gnagspår"
$ws.Range("Q60").Value = 557081
$ws.Range("R60").Value = 6710301
$ws.Range("Z60").Value = "11:20"
$ws.Range("AB60").Value = "11:20"

# Row 61
$ws.Range("A61").Value = 130671331
$ws.Range("B61").Value = 91808
$ws.Range("D61").Value = "NT"
$ws.Range("E61").Value = 1202
$ws.Range("F61").Value = "Ullticka"
$ws.Range("G61").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H61").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("M61").ClearContents()
$ws.Range("Q61").Value = 556976
$ws.Range("R61").Value = 6710393
$ws.Range("Z61").Value = "11:04"
$ws.Range("AB61").Value = "11:04"
